# Apply "Running IAM test cases" edit:
# - Rows 3-7 in column C change from "Y" to "N" (row 2 stays "Y")
# - Selection moves from D8 to C4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"

$ws.Activate()
$ws.Range("C4").Select()
